$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A values (rows 2-39) after the re-shuffle
$values = @{
    2  = "40-4-2"
    3  = "52-2-1"
    4  = "37-5-1"
    5  = "48-4-1"
    6  = "52-2-2"
    7  = "52-3-1"
    8  = "40-2-2"
    9  = "45-5-2"
    10 = "42-3-2"
    11 = "50-3-1"
    12 = "49-1-1"
    13 = "30-4-3"
    14 = "49-3-2"
    15 = "49-2-2"
    16 = "66-3-1"
    17 = "38-1-1"
    18 = "56-2-2"
    19 = "50-3-2"
    20 = "50-4-1"
    21 = "40-5-1"
    22 = "46-3-2"
    23 = "47-4-2"
    24 = "51-2-2"
    25 = "38-4-1"
    26 = "49-1-2"
    27 = "37-1-1"
    28 = "51-1-1"
    29 = "41-3-1"
    30 = "46-3-1"
    31 = "50-2-1"
    32 = "46-1-2"
    33 = "50-1-2"
    34 = "46-4-1"
    35 = "47-1-1"
    36 = "45-5-1"
    37 = "46-4-0"
    38 = "38-2-1"
    39 = "49-2-1"
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 1).Value = $values[$row]
}

# Row 3 gains new "unknown emotion in:" / "unknown author in:" entries in B and C
$ws.Range("B3").Value = "46-4-0"
$ws.Range("C3").Value = "46-4-0"
